# HotelResData.xlsx — give the first four hotel reservations (HR001-HR004)
# an earlier stay (04/02/2019 -> 06/02/2019) instead of the 10/02/2019 ->
# 12/02/2019 block shared by all eight rows. HR005-HR008 keep their
# original dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "04/02/2019"
$ws.Range("D2").Value = "06/02/2019"
$ws.Range("C3").Value = "04/02/2019"
$ws.Range("D3").Value = "06/02/2019"
$ws.Range("C4").Value = "04/02/2019"
$ws.Range("D4").Value = "06/02/2019"
$ws.Range("C5").Value = "04/02/2019"
$ws.Range("D5").Value = "06/02/2019"

# Leave the cursor where the author's save left it.
$ws.Range("B10").Select()
